# ============================================================================
# Edit script: add two new per-stage summary sheets (ariane_only_stg1,
# bsg_only_stg1), refresh the bsg_chip sheet layout/values, and update the
# view/selection state to match the post-edit workbook.
# ============================================================================

$wb = $excel.ActiveWorkbook

$sheet1    = $wb.Worksheets.Item(1)   # Sheet1
$arianeFix = $wb.Worksheets.Item(2)   # ariane_fix
$bsgChip   = $wb.Worksheets.Item(3)   # bsg_chip

# ----------------------------------------------------------------------------
# 1. bsg_chip: drop the "no_stg2" / "no_stg2_sort" columns (D:E), which
#    shifts the trailing "bsg_chip" (now "of_after_stg1") column left to D,
#    then refresh the "sort" column (C) values and the new D header/value.
# ----------------------------------------------------------------------------
$bsgChip.Range("D1:E1").EntireColumn.Delete() | Out-Null

$bsgChip.Range("D1").Value = "of_after_stg1"

$bsgChip.Range("C2").ClearContents() | Out-Null
$bsgChip.Range("C3").ClearContents() | Out-Null
$bsgChip.Range("C4").Value = 16399020.770132
$bsgChip.Range("C5").Value = 16233372
$bsgChip.Range("C6").Value = 33448088.539326001
$bsgChip.Range("C9").ClearContents() | Out-Null

$bsgChip.Columns.Item(1).ColumnWidth = 15.5
$bsgChip.Columns.Item(2).ColumnWidth = 12.75
$bsgChip.Columns.Item(3).ColumnWidth = 14.75

$bsgChip.Range("C9").Select() | Out-Null

# ----------------------------------------------------------------------------
# 2. ariane_fix: clear the old tab selection / cell selection.
# ----------------------------------------------------------------------------
$arianeFix.Range("B2:C2").Select() | Out-Null

# ----------------------------------------------------------------------------
# 3. New sheet "ariane_only_stg1" - stage-1-only stats for ariane.
# ----------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$arianeStg1 = $wb.Worksheets.Add($null, $last)
$arianeStg1.Name = "ariane_only_stg1"

$bsgChip.Range("B1:C1").Copy() | Out-Null
$arianeStg1.Range("B1").PasteSpecial(-4122) | Out-Null
$arianeStg1.Range("B1").Value = "baseline"
$arianeStg1.Range("C1").Value = "sort"

$bsgChip.Range("A2:C7").Copy() | Out-Null
$arianeStg1.Range("A2").PasteSpecial(-4122) | Out-Null

$arianeStg1.Range("A2").Value = "of_after_stg1"
$arianeStg1.Range("B2").Value = 26052
$arianeStg1.Range("C2").Value = 21487

$arianeStg1.Range("A3").Value = "wl"
$arianeStg1.Range("B3").Value = 3767357.2043110002
$arianeStg1.Range("C3").Value = 3767357.2043110002

$arianeStg1.Range("A4").Value = "via"
$arianeStg1.Range("B4").Value = 3301836
$arianeStg1.Range("C4").Value = 3297764

$arianeStg1.Range("A5").Value = "of"
$arianeStg1.Range("B5").Value = 10316379.26533
$arianeStg1.Range("C5").Value = 10258074.674814999

$arianeStg1.Range("A6").Value = "wns"
$arianeStg1.Range("B6").Value = -0.49
$arianeStg1.Range("C6").Value = -0.49

$arianeStg1.Range("A7").Value = "tns"
$arianeStg1.Range("B7").Value = -1210.6199999999999
$arianeStg1.Range("C7").Value = -1201.3699999999999

$arianeStg1.Columns.Item(1).ColumnWidth = 12.125
$arianeStg1.Columns.Item(2).ColumnWidth = 11.625

$arianeStg1.Range("E18").Select() | Out-Null

# ----------------------------------------------------------------------------
# 4. New sheet "bsg_only_stg1" - stage-1-only stats for bsg_chip.
# ----------------------------------------------------------------------------
$bsgStg1 = $wb.Worksheets.Add($null, $arianeStg1)
$bsgStg1.Name = "bsg_only_stg1"

$bsgChip.Range("B1:C1").Copy() | Out-Null
$bsgStg1.Range("B1").PasteSpecial(-4122) | Out-Null
$bsgStg1.Range("B1").Value = "baseline"
$bsgStg1.Range("C1").Value = "sort"

$bsgChip.Range("A2:C7").Copy() | Out-Null
$bsgStg1.Range("A2").PasteSpecial(-4122) | Out-Null

$bsgStg1.Range("A2").Value = "of_after_stg1"
$bsgStg1.Range("B2").Value = 110273
$bsgStg1.Range("C2").ClearContents() | Out-Null

$bsgStg1.Range("A3").Value = "wl"
$bsgStg1.Range("B3").Value = 21497469.125300001
$bsgStg1.Range("C3").Value = 21497469.125300001

$bsgStg1.Range("A4").Value = "via"
$bsgStg1.Range("B4").Value = 20866152
$bsgStg1.Range("C4").Value = 20851260

$bsgStg1.Range("A5").Value = "of"
$bsgStg1.Range("B5").Value = 34961494.361810997
$bsgStg1.Range("C5").Value = 34667571.438350998

$bsgStg1.Range("A6").Value = "wns"
$bsgStg1.Range("B6").Value = -0.42
$bsgStg1.Range("C6").Value = -0.42

$bsgStg1.Range("A7").Value = "tns"
$bsgStg1.Range("B7").Value = -8842.82
$bsgStg1.Range("C7").Value = -8767.32

$bsgStg1.Columns.Item(1).ColumnWidth = 12.125
$bsgStg1.Columns.Item(2).ColumnWidth = 12.75
$bsgStg1.Columns.Item(3).ColumnWidth = 12.75

$bsgStg1.Range("D14").Select() | Out-Null

# ----------------------------------------------------------------------------
# 5. Final active sheet/tab is the last one added (bsg_only_stg1).
# ----------------------------------------------------------------------------
$bsgStg1.Activate()
